$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") holds a count derived per-row. Regenerate the
# values (K instead of Strike#) for rows 2-12.
$kValues = @{
    2  = 4
    3  = 5
    4  = 6
    5  = 7
    6  = 6
    7  = 5
    8  = 2
    9  = 1
    10 = 1
    11 = 3
    12 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
